# "Add files via upload" / "ver 26/10/2023"
#
# Append four newly-collected survey rows (34-37) to the bottom of the
# "Hoja2" table, continuing the existing Fecha/Encuestadora/... columns.
# Formatting is carried down from the last existing row (33) so the new
# rows look exactly like the rest of the table (date format in A, the
# Encuestadora fill/center style in B, percentage columns in C:J), then
# the real values are written on top.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 33's formatting down across the 4 new rows first (mirrors
# selecting row 33, copying, and pasting formats onto 34:37 in the UI).
$ws.Range("A33:J33").Copy()
$ws.Range("A34:J37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 34: 2023-10-23, RCN
$ws.Cells.Item(34, 1).Value = 45222
$ws.Cells.Item(34, 2).Value = "RCN"
$ws.Cells.Item(34, 3).Value = 0.369
$ws.Cells.Item(34, 4).Value = 0.288
$ws.Cells.Item(34, 5).Value = 0.088
$ws.Cells.Item(34, 6).Value = 0.09
$ws.Cells.Item(34, 7).Value = 0.097
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 0.06800000000000006

# Row 35: 2023-10-24, RCN
$ws.Cells.Item(35, 1).Value = 45223
$ws.Cells.Item(35, 2).Value = "RCN"
$ws.Cells.Item(35, 3).Value = 0.372
$ws.Cells.Item(35, 4).Value = 0.29199999999999998
$ws.Cells.Item(35, 5).Value = 0.086
$ws.Cells.Item(35, 6).Value = 0.093
$ws.Cells.Item(35, 7).Value = 0.088
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 10).Value = 0.06900000000000017

# Row 36: 2023-10-25, RCN
$ws.Cells.Item(36, 1).Value = 45224
$ws.Cells.Item(36, 2).Value = "RCN"
$ws.Cells.Item(36, 3).Value = 0.36399999999999999
$ws.Cells.Item(36, 4).Value = 0.29899999999999999
$ws.Cells.Item(36, 5).Value = 0.09
$ws.Cells.Item(36, 6).Value = 0.09
$ws.Cells.Item(36, 7).Value = 0.088
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 10).Value = 0.06900000000000006

# Row 37: 2023-10-25, Invamer
$ws.Cells.Item(37, 1).Value = 45224
$ws.Cells.Item(37, 2).Value = "Invamer"
$ws.Cells.Item(37, 3).Value = 0.35499999999999998
$ws.Cells.Item(37, 4).Value = 0.276
$ws.Cells.Item(37, 5).Value = 0.111
$ws.Cells.Item(37, 6).Value = 0.091
$ws.Cells.Item(37, 7).Value = 0.10299999999999999
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 0.06400000000000006

# Leave the view/selection where the user ended up after entering the data.
$ws.Range("G34").Select()
